$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumPoliza (E2) changes from 04005238448 to 04104013020.
# The cell is formatted as Text with a quote-prefix (stored digits string),
# so assign the value with a leading apostrophe to keep it text and keep
# the existing "quote prefix" cell style instead of Excel creating a new one.
$ws.Range("E2").Value = "'04104013020"

# The active selection moved to F11 before the file was saved.
$ws.Range("F11").Select()
